# Weekly refresh of the "Perejil" (Vega Monumental Concepción) price series.
# A new week's observations (Primera/Segunda) are inserted right after the
# existing row 53 block, which pushes every later row down by two rows and
# naturally turns the former last pair of rows into the new last pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 54:55 - rows 54..231 shift down to 56..233.
$ws.Rows("54:55").Insert()

# Seed the new pair of rows with the same shape/values as the row pair that
# is now directly below them (the old row 54:55, now at 56:57), then update
# just the date to the new week being added.
$ws.Range("A54:R55").Value2 = $ws.Range("A56:R57").Value2
$ws.Range("D54:D55").Value2 = 45133
